# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (column E) and
# "Correspond Handback DateTime" (column H) timestamps on row 2
# (the 73bb9739-... / 11ed010e... file) of both the "zh-cn" and
# "de-de" report sheets to reflect a newly generated handback report.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E2").Value = "2016-03-18 12:34:04"
$zhcn.Range("H2").Value = "2016-03-18 12:34:26"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E2").Value = "2016-03-18 12:34:07"
$dede.Range("H2").Value = "2016-03-18 12:34:33"
